$wb = $excel.ActiveWorkbook

# --- tipo_persona (sheet1) ---
$ws1 = $wb.Worksheets.Item("tipo_persona")
$table = "tipo_persona"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws1.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws1.Range("D3").Formula = $f3

# --- tipo_persona_rel (sheet2) ---
$ws2 = $wb.Worksheets.Item("tipo_persona_rel")
$table = "tipo_persona_rel"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws2.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws2.Range("D3").Formula = $f3

# --- operacion_titulo (sheet3) ---
$ws3 = $wb.Worksheets.Item("operacion_titulo")
$table = "operacion_titulo"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws3.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws3.Range("D3").Formula = $f3
$f4 = '="insert into interno.' + $table + " values ('" + '" & A4 & "' + "','" + '" & B4 & "' + "');" + '"'
$ws3.Range("D4").Formula = $f4
$f5 = '="insert into interno.' + $table + " values ('" + '" & A5 & "' + "','" + '" & B5 & "' + "');" + '"'
$ws3.Range("D5").Formula = $f5

# --- operacion_titulo_rel (sheet4) ---
$ws4 = $wb.Worksheets.Item("operacion_titulo_rel")
$table = "operacion_titulo_rel"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws4.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws4.Range("D3:D5").Formula = $f3

# --- tipo_deudor (sheet5) ---
$ws5 = $wb.Worksheets.Item("tipo_deudor")
$table = "tipo_deudor"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws5.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws5.Range("D3").Formula = $f3

# --- tipo_deudor_rel (sheet6) ---
$ws6 = $wb.Worksheets.Item("tipo_deudor_rel")
$table = "tipo_deudor_rel"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws6.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws6.Range("D3").Formula = $f3

# --- tabla_banco_126 (sheet7) ---
$ws7 = $wb.Worksheets.Item("tabla_banco_126")
$table = "tabla_banco_126"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws7.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws7.Range("D3:D24").Formula = $f3

# --- tabla_banco_126_rel (sheet8) ---
$ws8 = $wb.Worksheets.Item("tabla_banco_126_rel")
$table = "tabla_banco_126_rel"
$f2 = '="insert into interno.' + $table + " values ('" + '" & A2 & "' + "','" + '" & B2 & "' + "');" + '"'
$ws8.Range("D2").Formula = $f2
$f3 = '="insert into interno.' + $table + " values ('" + '" & A3 & "' + "','" + '" & B3 & "' + "');" + '"'
$ws8.Range("D3:D24").Formula = $f3

# Row-height tweak observed on the banking-code sheet (two wrapped rows
# shrank from 36 to 24 points).
$ws7.Rows.Item(9).RowHeight = 24
$ws7.Rows.Item(17).RowHeight = 24

# Restore each sheet's own lingering selection, then finish on
# tabla_banco_126_rel / K9 as the active sheet+cell.
$ws1.Range("D2").Select()
$ws2.Range("D2").Select()
$ws3.Range("D2").Select()
$ws4.Range("D2").Select()
$ws5.Range("D2:D3").Select()
$ws6.Range("D2").Select()
$ws7.Range("D2").Select()
$ws8.Range("K9").Select()
